$d = $word.ActiveDocument
$find = $d.Content.Find
$find.Execute(", and report to jail on October 14, 2022, at 7:00 p.m", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
